$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.876.68"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "2.603.04"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.91"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.97"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +1.22%  "
$ws.Range("E9").Value = "  +2.10%  "
$ws.Range("E10").Value = "  +1.70%  "
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").Value = "3.058.89"
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("D14").Value = "60.895.35"
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.67"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("E16").Value = "  +0.96%  "
$ws.Range("D17").Value = "2.604.05"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.76"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "355.21"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.55%  "
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("E21").Value = "  +1.59%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.00"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.41%  "
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.167"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "2.718.60"
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.992"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.40"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.28"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +10.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.43"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.24%  "
$ws.Range("E33").Value = "  +2.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.19"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.20%  "
$ws.Range("E35").Value = "  +5.19%  "
$ws.Range("E36").Value = "  +1.04%  "
$ws.Range("E37").Value = "  +7.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.912"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +7.62%  "
$ws.Range("E39").Value = "  +1.25%  "
$ws.Range("E40").Value = "  +1.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "292.18"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.40%  "
$ws.Range("E43").Value = "  +1.77%  "
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.61"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.78%  "
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("E49").Value = "  +2.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.33"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("E51").Value = "  +8.05%  "
